$d = $word.ActiveDocument

# --- Step 1: simplify the "Meteo France" paragraph (merge split runs; same visible text) ---
$oldText1 = "#process SEVIRI data to create RGB images with dust mask (II method...Meteo France)"
$null = $d.Content.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $oldText1, 2)

# --- Step 2: delete the old, now-duplicated SEVIRI/MODIS block before inserting the new one ---
# (Doing the delete first, while the marker text is still unique, avoids ambiguity.)
$full = $d.Content.Text
$startMarker = "# SEVIRI RGB mask daytime"
$endMarker = "30 17 * * * /home/mariners/MODIS_AOD/MODIS_AOD_download_hdf_ocean.sh 1>/home/mariners/log/run_MODIS.log 2>/home/mariners/log/run_MODIS-error.log"
$startIdx = $full.IndexOf($startMarker)
$endIdx = $full.IndexOf($endMarker) + $endMarker.Length
$delRange = $d.Range($startIdx, $endIdx + 1)
$delRange.Delete()

# --- Step 3: insert the new block of paragraphs right after "#run contrab every 35th minute" ---
$full2 = $d.Content.Text
$anchor = "#run contrab every 35th minute"
$anchorIdx = $full2.IndexOf($anchor)
$insertAt = $anchorIdx + $anchor.Length
$insertPoint = $d.Range($insertAt, $insertAt)
$newBlock = "`r# SEVIRI RGB mask daytime`r*/45 * * * * /bin/sh /home/mariners/SEVIRI_DUST/execute_SEVIRI_mask.sh  1>/home/mariners/log/run.log 2>/home/mariners/`$`r`r# SEVIRI RGB mask nighttime`r*/50 * * * * /bin/sh /home/mariners/SEVIRI_DUST/execute_SEVIRI_mask_NIGHT.sh  1>/home/mariners/log/run_NIGHT.log 2>/ho`$`r`r# generate list of solar_zenith angle at the end of the day every 20 minutes`r*/20 * * * * /bin/sh /home/mariners/SEVIRI_DUST/execute_SEVIRI_solar_zenith.sh 1>/home/mariners/log/run_Zenith.log 2>/`$`r`r######################################################################################################################`$`r# remove SEVIRI .img files older than 7 days`r30 08 * * * bash /home/mariners/SEVIRI_DUST/scripts/remove_seviriData.bash 1>/home/mariners/log/del.log 2>/home/marine`$`r`r######################################################################################################################`$`r# crontab to download data from AOD MODIS-TERRA & AQUA (10km)`r30 18 * * * /home/mariners/MODIS_AOD/MODIS_AOD_download_hdf_ocean.sh 1>/home/mariners/log/run_MODIS.log 2>/home/marine`$`r`r######################################################################################################################`$`r# crontab to create DAILY Air Quality Indexes from MODIDS data`r45 18 * * * /home/mariners/MODIS_AOD/AQI_MODIS_AOD_hdf_ocean.sh 1>/home/mariners/log/run_AQI.log 2>/home/mariners/log/`$`r`r######################################################################################################################`$`r# generate DUST allerts during`r*/50 * * * * bash /home/mariners/SEVIRI_DUST/scripts/latest_MASK.sh 1>/home/mariners/log/allerts.log 2>/home/mariners/`$`r"
$insertPoint.InsertAfter($newBlock)

Write-Output ("final paragraph count=" + $d.Paragraphs.Count)
